$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")

# Row 16: change Method from GET to POST, and URL to the trimmed login endpoint
$ws.Range("A16").Value = "POST"
$ws.Range("B16").Value = "/api/users/login/"

# Row 18: add new Registeration endpoint row
$ws.Range("A18").Value = "POST"
$ws.Range("B18").Value = "/api/users/registeration"

# Update selection to reflect where the user last left off (C18)
$ws.Range("C18").Select()
